# country -> group, countries -> groupings
$wb = $excel.ActiveWorkbook

$params = $wb.Worksheets.Item("params")
$timeA  = $wb.Worksheets.Item("time_A")
$timeB  = $wb.Worksheets.Item("time_B")
$powerA = $wb.Worksheets.Item("power_A")

# Rename header cells from "country" to "group"
$params.Range("C1").Value = "group"
$timeA.Range("A1").Value = "group"
$timeB.Range("A1").Value = "group"
$powerA.Range("A1").Value = "group"

# Move the active selection: params becomes the active sheet with C1 selected,
# and time_A is no longer the tab-selected sheet (A1 stays selected there).
$null = $timeA.Range("A1").Select()
$null = $params.Activate()
$null = $params.Range("C1").Select()
